$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 3 (castor) content updates ---
$ws.Range("C3").Value = "data-raw/logo_castor.jpg"
$ws.Range("D3").Value = "Quasi disparu en Europe au début du 20ème siècle, la prise de mesure de protection du Castor d’Europe à partir de 1909 puis son classement en « espèce protégée » en 1968 a permis leur reconquête du territoire. L’ Ile-de-France est l’un des fronts de colonisation, induisant un suivi annuel et précis permettant de détecter la présence, suivre la distribution de l’espèce et la protéger en conséquence (notamment via l’interdiction de piégeage)."
$ws.Range("E3").Value = "Les objectifs du réseau ont été fixés par le ministère en charge de l’écologie :`n- assurer le suivi de la colonisation du castor sur le réseau hydrographique français ;`n- accompagnement sur la question des dommages;`n- vigilance sur l'arrivée du castor canadien"
$ws.Range("F3").Value = "Raportage Directive européenne Habitat Faune Flore`nRégulation des dispositifs de piègeage près des cours d'eau"
$ws.Range("G3").Value = "77, 91"
$ws.Range("H3").Value = "data-raw/lineaire_castor.gpkg"
$ws.Range("I3").Value = "Prospection de linéaires de cours d'eau à la recherche d'indices de présence"
$ws.Range("K3").Value = "Prospections préférentiellement hivernales avant la reprise de la végétation"
$ws.Range("L3").Value = "Animation nationale:`nPaul Hurel`nSuivi scientifique:`nYoann Bressan`nCouriel du réseau:`nreseau.castor@ofb.gouv.fr`nAnimation régionale:`nCédric Mondy"
$ws.Range("M3").Value = "Conseils départementaux (ENS)`nSyndicats de rivière"
$ws.Range("O3").Value = "environ 1/2 journée par prospection"
$ws.Range("Q3").Value = "4j"
$ws.Range("S3").Value = "Coordination`nRemontée des données au national"
$ws.Range("T3").Value = "Prospections`n(Constats de dommage)"
$ws.Range("U3").Value = "Formation Petit et Méso-Carnivores et Castor`n(Formation dommage)"
$ws.Range("V3").Value = "Recherche d'indices de présence (bois coupé, écorçage, hutte…) en prospection sur l'eau et à pied en berge"
$ws.Range("W3").Value = "- embarquation (ex. kayak)`n- gilet de sauvetage`n- jumelles`n- appareil photo`n- GPS"
$ws.Range("Y3").Value = "Remplissage des fiches terrains`nBancarisation régionale`nTransmission au national qui effectue une validation et consolidation nationale des données"
$ws.Range("Z3").Value = "https://carmen.carmencarto.fr/38/Castor.map"
$ws.Range("AA3").Value = "texte:Le réseau castor;lien:https://professionnels.ofb.fr/fr/reseau-castor"
$ws.Range("AB3").Value = "texte:Fiche espèce;lien:https://professionnels.ofb.fr/fr/doc-fiches-especes/castor-deurope-castor-fiber"
$ws.Range("AD3").Value = "texte: Site Alfresco du Réseau Castor IdF;lien:https://ged.ofb.fr/share/page/site/dridf-rseau-partenarial-castor/dashboard"
$ws.Range("AE3").Value = "texte: Protocole;lien:https://ged.ofb.fr/share/s/giB4EPFIRPmsQZiGFeYY0A"
$ws.Range("AF3").Value = "texte:Arrêtés piégeage (Serveur DR);lien:\\ad.intra\dfs\COMMUNS\REGIONS\IDF\DR\05_CONNAISSANCE\Castor\04_ArretesPiegeage"

# --- Column D width (now fits the longer description text) ---
$ws.Columns.Item(4).ColumnWidth = 36.3

# --- Row 3 height (shrinks now that text is reorganised) ---
$ws.Rows.Item(3).RowHeight = 195

# --- Update the frozen-pane view position / active selection ---
$ws.Range("AC3").Select()
